$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the old emoji markers used in column A (statut) to their replacement text,
# as described by the diff of the shared strings table.
$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📙" = "+3"
    "📗" = "✅"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Values that look like numbers ("-3", "+3") need a leading apostrophe so
# Excel stores them as text (shared string) instead of converting them to a
# numeric value (which would also silently drop the "+" sign).
$textForce = @{
    "-3" = $true
    "+3" = $true
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $newValue = $map[$current]
        if ($textForce.ContainsKey($newValue)) {
            $cell.Value2 = "'" + $newValue
        } else {
            $cell.Value2 = $newValue
        }
    }
}
